$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A width: 10.7109375 -> 11.7109375 (character width, as stored in OOXML).
# The COM ColumnWidth setter here quantizes to the nearest 1/6 of a character,
# so 10.833333 is the input that lands on the closest representable width
# (11.666666666666666) to the target 11.7109375.
$ws.Columns.Item(1).ColumnWidth = 10.833333333333332

# Updated data values for row 1
$ws.Range("A1").Value = 149.11057993586851
$ws.Range("B1").Value = 3.8486190834476908
$ws.Range("C1").Value = 2.4624918460534899
